$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 723 (the "teacher-less French" post). This shifts all
# subsequent rows (724:824) up by one, matching the new dimension A1:C823.
$ws.Rows.Item(723).Delete()
